$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 survey response: "wantUSCurrency" answer corrected from "Yes" to "No"
$ws.Range("D7").Value = "No"

# Update the active selection to reflect where the editor left off
$ws.Range("D16").Select() | Out-Null
